$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
    $ws.Range("H6").Value = 3509031.8
    $ws.Range("I6").Value = 4048817.5
    $ws.Range("K6").Value = 12146452.5
    $ws.Range("M6").Value = -12146340.5
    $ws.Range("H64").Value = 3818.348
    $ws.Range("J64").Value = 3891.389
    $ws.Range("L64").Value = 3891.389
    $ws.Range("N64").Value = -4387.389
    $ws.Range("H67").Value = 3818.348
    $ws.Range("J67").Value = 3891.389
    $ws.Range("L67").Value = 3891.389
    $ws.Range("N67").Value = -5607.389
    $ws.Range("H97").Value = 601.4286
    $ws.Range("J97").Value = 601.4286
    $ws.Range("L97").Value = 1804.2858
    $ws.Range("N97").Value = -2796.2858
    $ws.Range("H98").Value = 2260
    $ws.Range("I98").Value = 2150
    $ws.Range("K98").Value = 2150
    $ws.Range("M98").Value = -652
    $ws.Range("H100").Value = 1400
    $ws.Range("I100").Value = 0
    $ws.Range("J100").Value = 1400
    $ws.Range("K100").Value = 0
    $ws.Range("L100").Value = 1400
    $ws.Range("N100").Value = -2482
    $ws.Range("H101").Value = 466.75
    $ws.Range("I101").Value = 355.33334
    $ws.Range("J101").Value = 801
    $ws.Range("K101").Value = 1066.00002
    $ws.Range("L101").Value = 2403
    $ws.Range("M101").Value = 555.9999800000001
    $ws.Range("N101").Value = -5647
    $ws.Range("H122").Value = 2260
    $ws.Range("I122").Value = 2150
    $ws.Range("K122").Value = 6450
    $ws.Range("M122").Value = -4000
    $ws.Range("H125").Value = 2000
    $ws.Range("I125").Value = 2000
    $ws.Range("J125").Value = 2000
    $ws.Range("K125").Value = 18000
    $ws.Range("L125").Value = 18000
    $ws.Range("M125").Value = -15540
    $ws.Range("N125").Value = -22920
    $ws.Range("M100").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
    $ws.Range("H102").Value = 2902.5
    $ws.Range("I102").Value = 2902.5
    $ws.Range("K102").Value = 2902.5
    $ws.Range("M102").Value = -1280.5

$ws = $wb.Worksheets.Item("BSM")
    $ws.Range("H41").Value = 0
    $ws.Range("J41").Value = 0
    $ws.Range("L41").Value = 0
    $ws.Range("H48").Value = 0
    $ws.Range("J48").Value = 0
    $ws.Range("L48").Value = 0
    $ws.Range("H94").Value = 263.45
    $ws.Range("I94").Value = 241.8125
    $ws.Range("J94").Value = 350
    $ws.Range("K94").Value = 241.8125
    $ws.Range("L94").Value = 350
    $ws.Range("M94").Value = 209.1875
    $ws.Range("N94").Value = -1252
    $ws.Range("H99").Value = 2448.6875
    $ws.Range("I99").Value = 1026.5
    $ws.Range("J99").Value = 3302
    $ws.Range("K99").Value = 1026.5
    $ws.Range("L99").Value = 3302
    $ws.Range("M99").Value = 471.5
    $ws.Range("N99").Value = -6298
    $ws.Range("H105").Value = 3250.25
    $ws.Range("I105").Value = 2996.6667
    $ws.Range("K105").Value = 2996.6667
    $ws.Range("M105").Value = -1249.6667
    $ws.Range("H107").Value = 3101.5
    $ws.Range("I107").Value = 3174.9333
    $ws.Range("J107").Value = 2000
    $ws.Range("K107").Value = 3174.9333
    $ws.Range("L107").Value = 2000
    $ws.Range("M107").Value = -1254.9333
    $ws.Range("N107").Value = -5840
    $ws.Range("N41").ClearContents()
    $ws.Range("N48").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
    $ws.Range("H15").Value = 4504
    $ws.Range("I15").Value = 4008
    $ws.Range("K15").Value = 4008
    $ws.Range("M15").Value = -3838
    $ws.Range("H102").Value = 23000
    $ws.Range("J102").Value = 23000
    $ws.Range("L102").Value = 23000
    $ws.Range("N102").Value = -27868
    $ws.Range("H107").Value = 709.0769
    $ws.Range("I107").Value = 541.8570999999999
    $ws.Range("J107").Value = 904.1667
    $ws.Range("K107").Value = 541.8570999999999
    $ws.Range("L107").Value = 904.1667
    $ws.Range("M107").Value = 1378.1429
    $ws.Range("N107").Value = -4744.1667
    $ws.Range("H141").Value = 30000
    $ws.Range("I141").Value = 30000
    $ws.Range("J141").Value = 0
    $ws.Range("K141").Value = 30000
    $ws.Range("L141").Value = 0
    $ws.Range("M141").Value = -24820
    $ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
    $ws.Range("H7").Value = 975
    $ws.Range("I7").Value = 950
    $ws.Range("J7").Value = 1000
    $ws.Range("K7").Value = 2850
    $ws.Range("L7").Value = 3000
    $ws.Range("M7").Value = -2738
    $ws.Range("N7").Value = -3224
    $ws.Range("H68").Value = 0
    $ws.Range("I68").Value = 0
    $ws.Range("J68").Value = 0
    $ws.Range("K68").Value = 0
    $ws.Range("L68").Value = 0
    $ws.Range("H71").Value = 0
    $ws.Range("I71").Value = 0
    $ws.Range("J71").Value = 0
    $ws.Range("K71").Value = 0
    $ws.Range("L71").Value = 0
    $ws.Range("H80").Value = 0
    $ws.Range("I80").Value = 0
    $ws.Range("J80").Value = 0
    $ws.Range("K80").Value = 0
    $ws.Range("L80").Value = 0
    $ws.Range("H83").Value = 0
    $ws.Range("I83").Value = 0
    $ws.Range("J83").Value = 0
    $ws.Range("K83").Value = 0
    $ws.Range("L83").Value = 0
    $ws.Range("H92").Value = 667027.6
    $ws.Range("I92").Value = 2000298.2
    $ws.Range("J92").Value = 392.3
    $ws.Range("K92").Value = 6000894.6
    $ws.Range("L92").Value = 1176.9
    $ws.Range("M92").Value = -5999646.6
    $ws.Range("N92").Value = -3672.9
    $ws.Range("H131").Value = 724.5161000000001
    $ws.Range("J131").Value = 949.4595
    $ws.Range("L131").Value = 2848.3785
    $ws.Range("N131").Value = -12928.3785
    $ws.Range("M68").ClearContents()
    $ws.Range("N68").ClearContents()
    $ws.Range("M71").ClearContents()
    $ws.Range("N71").ClearContents()
    $ws.Range("M80").ClearContents()
    $ws.Range("N80").ClearContents()
    $ws.Range("M83").ClearContents()
    $ws.Range("N83").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
    $ws.Range("H3").Value = 10176.667
    $ws.Range("I3").Value = 10176.667
    $ws.Range("J3").Value = 0
    $ws.Range("K3").Value = 10176.667
    $ws.Range("L3").Value = 0
    $ws.Range("M3").Value = -10060.667
    $ws.Range("H7").Value = 3000
    $ws.Range("J7").Value = 3000
    $ws.Range("L7").Value = 3000
    $ws.Range("N7").Value = -3224
    $ws.Range("H8").Value = 3000
    $ws.Range("J8").Value = 3000
    $ws.Range("L8").Value = 3000
    $ws.Range("N8").Value = -3278
    $ws.Range("H14").Value = 14334000
    $ws.Range("I14").Value = 14334000
    $ws.Range("J14").Value = 0
    $ws.Range("K14").Value = 14334000
    $ws.Range("L14").Value = 0
    $ws.Range("M14").Value = -14333832
    $ws.Range("H80").Value = 2977.1924
    $ws.Range("I80").Value = 2835
    $ws.Range("J80").Value = 3245.7778
    $ws.Range("K80").Value = 2835
    $ws.Range("L80").Value = 3245.7778
    $ws.Range("M80").Value = -1837
    $ws.Range("N80").Value = -5241.7778
    $ws.Range("H83").Value = 2977.1924
    $ws.Range("I83").Value = 2835
    $ws.Range("J83").Value = 3245.7778
    $ws.Range("K83").Value = 14175
    $ws.Range("L83").Value = 16228.889
    $ws.Range("M83").Value = -9183
    $ws.Range("N83").Value = -26212.889
    $ws.Range("H126").Value = 3144.2307
    $ws.Range("I126").Value = 2377
    $ws.Range("J126").Value = 3801.8572
    $ws.Range("K126").Value = 7131
    $ws.Range("L126").Value = 11405.5716
    $ws.Range("M126").Value = -4661
    $ws.Range("N126").Value = -16345.5716
    $ws.Range("H132").Value = 5810.9443
    $ws.Range("I132").Value = 8867.143
    $ws.Range("J132").Value = 3866.0908
    $ws.Range("K132").Value = 26601.429
    $ws.Range("L132").Value = 11598.2724
    $ws.Range("M132").Value = -24071.429
    $ws.Range("N132").Value = -16658.2724
    $ws.Range("N3").ClearContents()
    $ws.Range("N14").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
    $ws.Range("H40").Value = 5364.706
    $ws.Range("I40").Value = 5811.1113
    $ws.Range("K40").Value = 5811.1113
    $ws.Range("M40").Value = -5675.1113
    $ws.Range("H61").Value = 8050.294
    $ws.Range("I61").Value = 10446.25
    $ws.Range("J61").Value = 2300
    $ws.Range("K61").Value = 10446.25
    $ws.Range("L61").Value = 2300
    $ws.Range("M61").Value = -10244.25
    $ws.Range("N61").Value = -2704
    $ws.Range("H82").Value = 1734
    $ws.Range("I82").Value = 1734
    $ws.Range("J82").Value = 0
    $ws.Range("K82").Value = 1734
    $ws.Range("L82").Value = 0
    $ws.Range("M82").Value = -1373
    $ws.Range("H85").Value = 1734
    $ws.Range("I85").Value = 1734
    $ws.Range("J85").Value = 0
    $ws.Range("K85").Value = 1734
    $ws.Range("L85").Value = 0
    $ws.Range("M85").Value = -486
    $ws.Range("H87").Value = 35866.668
    $ws.Range("J87").Value = 35866.668
    $ws.Range("L87").Value = 35866.668
    $ws.Range("N87").Value = -38112.668
    $ws.Range("H90").Value = 35866.668
    $ws.Range("J90").Value = 35866.668
    $ws.Range("L90").Value = 107600.004
    $ws.Range("N90").Value = -118832.004
    $ws.Range("H113").Value = 8050.294
    $ws.Range("I113").Value = 10446.25
    $ws.Range("J113").Value = 2300
    $ws.Range("K113").Value = 10446.25
    $ws.Range("L113").Value = 2300
    $ws.Range("M113").Value = -8276.25
    $ws.Range("N113").Value = -6640
    $ws.Range("H133").Value = 50576.47
    $ws.Range("J133").Value = 50576.47
    $ws.Range("L133").Value = 50576.47
    $ws.Range("N133").Value = -55636.47
    $ws.Range("N82").ClearContents()
    $ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
    $ws.Range("H46").Value = 0
    $ws.Range("J46").Value = 0
    $ws.Range("L46").Value = 0
    $ws.Range("H114").Value = 42199
    $ws.Range("J114").Value = 42199
    $ws.Range("L114").Value = 42199
    $ws.Range("N114").Value = -50877
    $ws.Range("H134").Value = 0
    $ws.Range("J134").Value = 0
    $ws.Range("L134").Value = 0
    $ws.Range("H140").Value = 62279.668
    $ws.Range("J140").Value = 62279.668
    $ws.Range("L140").Value = 62279.668
    $ws.Range("N140").Value = -72639.66800000001
    $ws.Range("N46").ClearContents()
    $ws.Range("N134").ClearContents()

Write-Output "Updated cells; sets and clears applied."